$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original cell styles for the Price column, force text number
# format while writing (so numeric-looking strings like "1.00" or "72.20"
# are not silently reinterpreted as floating point numbers), then restore
# the original style so no stray style index is left on the cells.
$priceRange = $ws.Range("D2:D51")
$priceStyle = $priceRange.Style
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "72.136.61"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "2.714.71"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "600.31"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").Value = "176.77"
$ws.Range("E6").Value = "  -1.45%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "0.524"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").Value = "2.712.94"
$ws.Range("E9").Value = "  +3.22%  "
$ws.Range("D10").Value = "0.170"
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("E11").Value = "  +2.65%  "
$ws.Range("D12").Value = "0.356"
$ws.Range("E12").Value = "  +2.22%  "
$ws.Range("D13").Value = "5.03"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").Value = "3.211.47"
$ws.Range("E14").Value = "  +2.40%  "
$ws.Range("D15").Value = "0.0000186"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").Value = "72.064.38"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "26.38"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "2.707.43"
$ws.Range("E18").Value = "  +2.86%  "
$ws.Range("D19").Value = "12.21"
$ws.Range("E19").Value = "  +6.59%  "
$ws.Range("D20").Value = "8.11"
$ws.Range("E20").Value = "  +3.63%  "
$ws.Range("D21").Value = "371.90"
$ws.Range("E21").Value = "  -2.82%  "
$ws.Range("D22").Value = "4.19"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").Value = "2.03"
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").Value = "72.20"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "4.36"
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("D27").Value = "9.91"
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("D28").Value = "2.852.34"
$ws.Range("E28").Value = "  +2.97%  "
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("D30").Value = "0.0₃0989"
$ws.Range("E30").Value = "  +2.61%  "
$ws.Range("D31").Value = "8.13"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("D32").Value = "512.21"
$ws.Range("E32").Value = "  -7.04%  "
$ws.Range("D33").Value = "1.31"
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "164.26"
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").Value = "19.76"
$ws.Range("E37").Value = "  +2.62%  "
$ws.Range("D38").Value = "19.08"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "1.39"
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("E40").Value = "  -4.60%  "
$ws.Range("D41").Value = "1.82"
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "5.08"
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "2.61"
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("D45").Value = "0.335"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("D46").Value = "157.78"
$ws.Range("E46").Value = "  +3.30%  "
$ws.Range("D47").Value = "39.39"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("D48").Value = "3.77"
$ws.Range("E48").Value = "  +3.36%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "0.563"
$ws.Range("E49").Value = "  +5.09%  "
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").Value = "1.78"
$ws.Range("E50").Value = "  +5.20%  "
$ws.Range("D51").Value = "0.0769"
$ws.Range("E51").Value = "  +1.39%  "

$priceRange.Style = $priceStyle
